# Adds a "2022-Q4" quarter to the workbook:
#   - inserts a new summary row at the top of the "总计" sheet
#   - inserts a new "2022-Q4" worksheet (fund holdings detail) before "2022-Q3"
#
# Helper: write a value that LOOKS like a number but must be stored as TEXT
# (mirrors how the existing fund-detail sheets store things like "38.21" /
# "001694" as plain text, not numbers). Using NumberFormat="@" forces Excel
# to keep the literal text instead of auto-converting it to a number; the
# ClearFormats() afterwards drops the "@" text format again so the cell ends
# up with no style index at all, matching the rest of the sheet.
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for 2022-Q4
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 19
$total.Cells.Item(2, 4).Value = 5.93

# Re-apply the column-A style (bold/centered/bordered) that the other rows
# in column A use, by copying the format from the row right below it.
$total.Cells.Item(3, 1).Copy()
$total.Cells.Item(2, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with fund-holdings detail
# ---------------------------------------------------------------------
# Duplicate the biggest existing quarter sheet so the new sheet inherits the
# exact same sheet-level formatting (outline props, page margins, header
# row style, column-A style, ...), then trim/overwrite its contents.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4src = $wb.Worksheets.Item("2021-Q4")
$q4src.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template sheet (2021-Q4) has 45 data rows (2..46); our new data only
# needs 19 (2..20) so drop the extra rows, shifting cells up.
$q4.Range("A21:H46").Delete(-4162)

# Row 2..20 data: [code, name, size, position, pct, value, rank]
$rows = @(
    @("001694", "华安沪港深外延增长混合A", "38.21", "94.15", "4.30", "1.6430", 3),
    @("006879", "华安智能生活混合A", "28.87", "92.19", "4.32", "1.2472", 3),
    @("007460", "华安成长创新混合A", "16.41", "93.15", "4.61", "0.7565", 3),
    @("013621", "华安智能生活混合C", "11.19", "92.19", "4.32", "0.4834", 3),
    @("002621", "中欧消费主题股票A", "10.21", "94.00", "4.26", "0.4349", 3),
    @("007126", "博道远航混合A", "6.63", "94.00", "4.26", "0.2824", 3),
    @("002697", "中欧消费主题股票C", "5.42", "94.00", "4.26", "0.2309", 3),
    @("014754", "华安景气优选混合A", "5.19", "92.65", "4.31", "0.2237", 3),
    @("007127", "博道远航混合C", "4.92", "94.00", "4.26", "0.2096", 3),
    @("008405", "华泰紫金泰盈混合C", "3.21", "91.44", "4.26", "0.1367", 8),
    @("860038", "光大阳光智造混合B", "2.61", "92.46", "4.35", "0.1135", 2),
    @("014755", "华安景气优选混合C", "1.48", "92.65", "4.31", "0.0638", 3),
    @("016099", "华安成长创新混合C", "0.47", "93.15", "4.61", "0.0217", 3),
    @("010124", "兴银景气优选混合A", "0.45", "83.23", "4.54", "0.0204", 8),
    @("008404", "华泰紫金泰盈混合A", "0.48", "91.44", "4.26", "0.0204", 8),
    @("010125", "兴银景气优选混合C", "0.35", "83.23", "4.54", "0.0159", 8),
    @("860039", "光大阳光智造混合C", "0.26", "92.46", "4.35", "0.0113", 2),
    @("860018", "光大阳光智造混合A", "0.18", "92.46", "4.35", "0.0078", 2),
    @("014972", "华安沪港深外延增长混合C", "0.13", "94.15", "4.30", "0.0056", 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $i

    Set-TextValue $q4.Cells.Item($r, 2) $data[0]
    Set-TextValue $q4.Cells.Item($r, 3) $data[1]
    Set-TextValue $q4.Cells.Item($r, 4) $data[2]
    Set-TextValue $q4.Cells.Item($r, 5) $data[3]
    Set-TextValue $q4.Cells.Item($r, 6) $data[4]
    Set-TextValue $q4.Cells.Item($r, 7) $data[5]

    $q4.Cells.Item($r, 8).Value = $data[6]
}
